$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the C13 (100nF) row - it is no longer necessary (near CP2102).
# This is row 23 on the BOM sheet. Deleting the whole row shifts
# everything below it up by one row.
$ws.Rows.Item(23).Delete()

$ws.Rows.Item(23).Select()
